# Generate Report for Handoff
# Updates the "b.md" rows across the Overview, zh-cn and de-de sheets to
# reflect that a new handoff package was generated for b.md.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 3 is the b.md file.
# Columns: A=File Name, B=Path And Name, C=Extension, D=Publish URL,
#          E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-09-04 04:42:06"

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is the b.md file.
# Columns: A=Source File Name, B=File Extension, C=Status, D=Source Path,
#          E=Priority, F=Content Duplicate, G=Latest Handoff File,
#          H=Latest Handoff Datetime, I=Latest Target File,
#          J=Latest Handback File, K=Latest Handback DateTime,
#          L=Reference Tokens, M=To be localized, N=Dependency From,
#          O=Has metadata, P=Error Detail
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-04 04:41:58"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48b5b7859638ca0fa34dd0aad7be29316255fac9/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/accd407628cecd8e2ea4daeb875f8946d4279449/e2e/b.md."
# Widen the Error Detail column (P) to fit the new message - target width 40
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# de-de sheet: row 3 is the b.md file. Same column layout as zh-cn.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-04 04:42:06"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48b5b7859638ca0fa34dd0aad7be29316255fac9/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/accd407628cecd8e2ea4daeb875f8946d4279449/e2e/b.md."
# Widen the Error Detail column (P) to fit the new message - target width 40
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
